$d = $word.ActiveDocument

# Change 1: "This Programs goal is to provide that valid header data to the User. "
#        -> "The purpose is to check if a Website has valid header data. "
$d.Content.Find.Execute("This Programs goal is to provide that valid header data to the User. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The purpose is to check if a Website has valid header data. ", 2)

# Change 2: " (website), and provider the webpage headers to the user."
#        -> " and get headers data from the site."
$d.Content.Find.Execute(" (website), and provider the webpage headers to the user.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " and get headers data from the site.", 2)
